$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Remove the "What team are you on?" column (table col 7 / sheet col G) ---
# The host doesn't support ListColumn.Delete, so shift the remaining header
# labels left within the table bounds, drop the now-duplicate trailing
# column, then resize the table so it picks the (now correct) header names
# back up off the sheet.
$ws.Range("G1").Value = "What days are you coming in?"
$ws.Range("H1").Value = "Choose a desk from the list below:"
$ws.Range("I1").EntireColumn.Delete()

# --- Drop the second response row (row 3) entirely ---
$ws.Range("A3").EntireRow.Delete()

# --- Resize the table to the new A1:H2 bounds / 8 columns ---
$lo.Resize($ws.Range("A1:H2"))

# --- Replace the single remaining response row with the new submission ---
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 45828.6792939815
$ws.Range("C2").Value = 45828.6794675926
$ws.Range("D2").Value = "slone@hoopp.com"
$ws.Range("E2").Value = "Shaheer Lone"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Monday;Wednesday;Thursday;"
$ws.Range("H2").Value = "16-W529"
